$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 91
$ws.Range("B91").Value = 6782568
$ws.Range("F91").Value = "Sporting San Jose"
$ws.Range("G91").Value = "AD Guanacasteca"
$ws.Range("H91").Value = 1
$ws.Range("I91").Value = 1
$ws.Range("J91").Value = "D"
$ws.Range("K91").Value = 1.909
$ws.Range("L91").Value = 3.6
$ws.Range("M91").Value = 3.3
$ws.Range("N91").Value = 2
$ws.Range("O91").Value = 3.6
$ws.Range("P91").Value = 3.1
$ws.Range("Q91").Value = -0.5
$ws.Range("R91").Value = 2
$ws.Range("S91").Value = 1.8
$ws.Range("T91").Value = 2.5
$ws.Range("U91").Value = 1.825
$ws.Range("V91").Value = 1.975
$ws.Range("W91").Value = -1
$ws.Range("X91").Value = 2.6
$ws.Range("Y91").Value = -1
$ws.Range("Z91").Value = -1
$ws.Range("AA91").Value = 0.8
$ws.Range("AB91").Value = -1
$ws.Range("AC91").Value = 0.9750000000000001

# Row 92
$ws.Range("B92").Value = 6782566
$ws.Range("F92").Value = "Cartagines"
$ws.Range("G92").Value = "Deportivo Saprissa"
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 4
$ws.Range("J92").Value = "A"
$ws.Range("K92").Value = 3.2
$ws.Range("L92").Value = 3.4
$ws.Range("M92").Value = 2
$ws.Range("N92").Value = 2.9
$ws.Range("O92").Value = 3.5
$ws.Range("P92").Value = 2.15
$ws.Range("Q92").Value = 0.25
$ws.Range("R92").Value = 1.875
$ws.Range("S92").Value = 1.925
$ws.Range("T92").Value = 3
$ws.Range("U92").Value = 1.975
$ws.Range("V92").Value = 1.825
$ws.Range("W92").Value = -1
$ws.Range("X92").Value = -1
$ws.Range("Y92").Value = 1.15
$ws.Range("Z92").Value = -1
$ws.Range("AA92").Value = 0.925
$ws.Range("AB92").Value = 0.9750000000000001
$ws.Range("AC92").Value = -1

# Row 129
$ws.Range("B129").Value = 6782595
$ws.Range("F129").Value = "Herediano"
$ws.Range("G129").Value = "Sporting San Jose"
$ws.Range("H129").Value = 3
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = "H"
$ws.Range("K129").Value = 1.4
$ws.Range("L129").Value = 4.75
$ws.Range("M129").Value = 7
$ws.Range("N129").Value = 1.363
$ws.Range("O129").Value = 4.75
$ws.Range("P129").Value = 8.5
$ws.Range("Q129").Value = -1.25
$ws.Range("R129").Value = 1.8
$ws.Range("S129").Value = 2
$ws.Range("T129").Value = 3
$ws.Range("U129").Value = 1.95
$ws.Range("V129").Value = 1.85
$ws.Range("W129").Value = 0.363
$ws.Range("X129").Value = -1
$ws.Range("Y129").Value = -1
$ws.Range("Z129").Value = 0.8
$ws.Range("AA129").Value = -1
$ws.Range("AB129").Value = 0
$ws.Range("AC129").Value = 0

# Row 130
$ws.Range("B130").Value = 6782598
$ws.Range("F130").Value = "Municipal Perez Zeledon"
$ws.Range("G130").Value = "Cartagines"
$ws.Range("H130").Value = 1
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = "H"
$ws.Range("K130").Value = 4.5
$ws.Range("L130").Value = 3.75
$ws.Range("M130").Value = 1.615
$ws.Range("N130").Value = 3.4
$ws.Range("O130").Value = 3.4
$ws.Range("P130").Value = 1.85
$ws.Range("Q130").Value = 0.5
$ws.Range("R130").Value = 1.8
$ws.Range("S130").Value = 2
$ws.Range("T130").Value = 2.75
$ws.Range("U130").Value = 1.9
$ws.Range("V130").Value = 1.9
$ws.Range("W130").Value = 2.4
$ws.Range("X130").Value = -1
$ws.Range("Y130").Value = -1
$ws.Range("Z130").Value = 0.8
$ws.Range("AA130").Value = -1
$ws.Range("AB130").Value = -1
$ws.Range("AC130").Value = 0.8999999999999999

# Row 131
$ws.Range("B131").Value = 6782596
$ws.Range("F131").Value = "Alajuelense"
$ws.Range("G131").Value = "AD Guanacasteca"
$ws.Range("H131").Value = 3
$ws.Range("I131").Value = 4
$ws.Range("J131").Value = "A"
$ws.Range("K131").Value = 1.363
$ws.Range("L131").Value = 4.75
$ws.Range("M131").Value = 8
$ws.Range("N131").Value = 1.444
$ws.Range("O131").Value = 4.333
$ws.Range("P131").Value = 7
$ws.Range("Q131").Value = -1.25
$ws.Range("R131").Value = 1.975
$ws.Range("S131").Value = 1.825
$ws.Range("T131").Value = 2.75
$ws.Range("U131").Value = 1.775
$ws.Range("V131").Value = 2.025
$ws.Range("W131").Value = -1
$ws.Range("X131").Value = -1
$ws.Range("Y131").Value = 6
$ws.Range("Z131").Value = -1
$ws.Range("AA131").Value = 0.825
$ws.Range("AB131").Value = 0.7749999999999999
$ws.Range("AC131").Value = -1

# Row 192
$ws.Range("B192").Value = 7623919
$ws.Range("F192").Value = "Municipal Liberia"
$ws.Range("G192").Value = "Sporting San Jose"
$ws.Range("H192").Value = 2
$ws.Range("I192").Value = 0
$ws.Range("J192").Value = "H"
$ws.Range("K192").Value = 1.75
$ws.Range("L192").Value = 3.6
$ws.Range("M192").Value = 3.8
$ws.Range("N192").Value = 1.8
$ws.Range("O192").Value = 3.6
$ws.Range("P192").Value = 3.6
$ws.Range("Q192").Value = -0.5
$ws.Range("R192").Value = 1.9
$ws.Range("S192").Value = 1.9
$ws.Range("T192").Value = 2.75
$ws.Range("U192").Value = 2
$ws.Range("V192").Value = 1.8
$ws.Range("W192").Value = 0.8
$ws.Range("X192").Value = -1
$ws.Range("Y192").Value = -1
$ws.Range("Z192").Value = 0.8999999999999999
$ws.Range("AA192").Value = -1
$ws.Range("AB192").Value = -1
$ws.Range("AC192").Value = 0.8

# Row 193
$ws.Range("B193").Value = 7623916
$ws.Range("F193").Value = "Santos de Gupiles"
$ws.Range("G193").Value = "AD Grecia"
$ws.Range("H193").Value = 0
$ws.Range("I193").Value = 2
$ws.Range("J193").Value = "A"
$ws.Range("K193").Value = 2.05
$ws.Range("L193").Value = 3.3
$ws.Range("M193").Value = 3.2
$ws.Range("N193").Value = 1.909
$ws.Range("O193").Value = 3.4
$ws.Range("P193").Value = 3.6
$ws.Range("Q193").Value = -0.5
$ws.Range("R193").Value = 1.95
$ws.Range("S193").Value = 1.85
$ws.Range("T193").Value = 2.5
$ws.Range("U193").Value = 1.85
$ws.Range("V193").Value = 1.95
$ws.Range("W193").Value = -1
$ws.Range("X193").Value = -1
$ws.Range("Y193").Value = 2.6
$ws.Range("Z193").Value = -1
$ws.Range("AA193").Value = 0.8500000000000001
$ws.Range("AB193").Value = -1
$ws.Range("AC193").Value = 0.95

# Row 205
$ws.Range("H205").Value = 3
$ws.Range("I205").Value = 1
$ws.Range("J205").Value = "H"
$ws.Range("N205").Value = 1.5
$ws.Range("O205").Value = 4.2
$ws.Range("P205").Value = 5.75
$ws.Range("Q205").Value = -1
$ws.Range("R205").Value = 1.8
$ws.Range("S205").Value = 2
$ws.Range("U205").Value = 1.775
$ws.Range("V205").Value = 2.025
$ws.Range("W205").Value = 0.5
$ws.Range("X205").Value = -1
$ws.Range("Y205").Value = -1
$ws.Range("Z205").Value = 0.8
$ws.Range("AA205").Value = -1
$ws.Range("AB205").Value = 0.7749999999999999
$ws.Range("AC205").Value = -1

# Row 206
$ws.Range("N206").Value = 1.8
$ws.Range("O206").Value = 3.5
$ws.Range("P206").Value = 4.75
$ws.Range("R206").Value = 2
$ws.Range("S206").Value = 1.8
$ws.Range("T206").Value = 2.5

# Row 207
$ws.Range("N207").Value = 1.285
$ws.Range("O207").Value = 5
$ws.Range("P207").Value = 11
$ws.Range("Q207").Value = -1.5
$ws.Range("R207").Value = 1.8
$ws.Range("S207").Value = 2
$ws.Range("T207").Value = 2.75
$ws.Range("U207").Value = 1.975
$ws.Range("V207").Value = 1.825

# Row 208
$ws.Range("N208").Value = 1.5
$ws.Range("O208").Value = 4
$ws.Range("P208").Value = 6.5
$ws.Range("Q208").Value = -1
$ws.Range("R208").Value = 1.85
$ws.Range("S208").Value = 1.95
$ws.Range("U208").Value = 1.85
$ws.Range("V208").Value = 1.95

# Row 209
$ws.Range("U209").Value = 2
$ws.Range("V209").Value = 1.8

# Row 210
$ws.Range("N210").Value = 5.75
$ws.Range("O210").Value = 4.5
$ws.Range("P210").Value = 1.5
$ws.Range("R210").Value = 2.025
$ws.Range("S210").Value = 1.775
$ws.Range("T210").Value = 2.75
$ws.Range("U210").Value = 1.825
$ws.Range("V210").Value = 1.975
